# This script applies a data correction to the "Germany Oberliga Hamburg" sheet.
# A handful of match rows had their data rows out of order / mismatched; this
# re-aligns the row data (everything except column A, the running index) by
# rotating values among the affected rows so each row's stats correspond to
# the correct match id / teams / odds.
#
# Column A (row index) and column D (match date, identical within each group)
# stay untouched; columns B through AD are moved between rows as described.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($row) {
    $rng = $ws.Range("B" + $row + ":AD" + $row)
    return $rng.Value2
}

function Set-RowValues($row, $values) {
    $rng = $ws.Range("B" + $row + ":AD" + $row)
    $rng.Value = $values
}

# Each group lists the rows involved, and for every row in the group it
# specifies which row's original data should end up there.
$groups = @(
    @{ Rows = @(20, 21, 22);      Source = @{ 20 = 22;  21 = 20;  22 = 21 } },
    @{ Rows = @(188, 189);        Source = @{ 188 = 189; 189 = 188 } },
    @{ Rows = @(205, 206);        Source = @{ 205 = 206; 206 = 205 } },
    @{ Rows = @(214, 215, 216);   Source = @{ 214 = 216; 215 = 214; 216 = 215 } },
    @{ Rows = @(222, 223);        Source = @{ 222 = 223; 223 = 222 } }
)

foreach ($group in $groups) {
    # Snapshot the current (pre-edit) values of every row in the group first,
    # so writes to one row never affect the data read for another.
    $snapshots = @{}
    foreach ($row in $group.Rows) {
        $snapshots[$row] = Get-RowValues $row
    }

    foreach ($row in $group.Rows) {
        $srcRow = $group.Source[$row]
        Set-RowValues $row $snapshots[$srcRow]
    }
}
